$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A39").Value = 45986
$ws.Range("B39").Value = 2025
$ws.Range("C39").Value = 0.21940175159154141
$ws.Range("D39").Value = 2026
$ws.Range("E39").Value = -0.18831859814396609

$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)

$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B39:E39").Select()
